# Weekly update: insert a new week's worth of records (2 rows) for
# "Vega Monumental Concepción - Coliflor" at the top of the data block,
# pushing the existing history down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows starting at row 197; everything currently at
# row 197 onward (through the old last row, 306) shifts down to 199..308.
$ws.Rows.Item(197).Resize(2).Insert()

# ---- New row 197: "Primera" quality record for the new date ----
$ws.Cells.Item(197, 1).Value = 11
$ws.Cells.Item(197, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(197, 3).Value = "Bíobío"
$ws.Cells.Item(197, 4).Value2 = 44813
$ws.Cells.Item(197, 5).Value = 8
$ws.Cells.Item(197, 6).Value = 100112008
$ws.Cells.Item(197, 7).Value = "Coliflor"
$ws.Cells.Item(197, 8).Value = "Sin especificar"
$ws.Cells.Item(197, 9).Value = "Primera"
$ws.Cells.Item(197, 10).Value = 2000
$ws.Cells.Item(197, 11).Value = 1200
$ws.Cells.Item(197, 12).Value = 1300
$ws.Cells.Item(197, 13).Value = 1250
$ws.Cells.Item(197, 14).Value = "$/unidad"
$ws.Cells.Item(197, 15).Value = "Región Metropolitana"
$ws.Cells.Item(197, 16).Value = 1250
$ws.Cells.Item(197, 17).Value = 1
$ws.Cells.Item(197, 18).Value = "Hortaliza"

# ---- New row 198: "Segunda" quality record for the new date ----
$ws.Cells.Item(198, 1).Value = 11
$ws.Cells.Item(198, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(198, 3).Value = "Bíobío"
$ws.Cells.Item(198, 4).Value2 = 44813
$ws.Cells.Item(198, 5).Value = 8
$ws.Cells.Item(198, 6).Value = 100112008
$ws.Cells.Item(198, 7).Value = "Coliflor"
$ws.Cells.Item(198, 8).Value = "Sin especificar"
$ws.Cells.Item(198, 9).Value = "Segunda"
$ws.Cells.Item(198, 10).Value = 1000
$ws.Cells.Item(198, 11).Value = 900
$ws.Cells.Item(198, 12).Value = 900
$ws.Cells.Item(198, 13).Value = 900
$ws.Cells.Item(198, 14).Value = "$/unidad"
$ws.Cells.Item(198, 15).Value = "Región Metropolitana"
$ws.Cells.Item(198, 16).Value = 900
$ws.Cells.Item(198, 17).Value = 1
$ws.Cells.Item(198, 18).Value = "Hortaliza"

# Keep the date column's existing number format on the two new rows
# (matches the "s=2" style used by every other Fecha cell in column D).
$ws.Range("D197:D198").NumberFormat = $ws.Range("D199").NumberFormat()
